$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rows whose "Observaciones" (column G) flag flips from "Actualizado " to
# "Actualizado *" (note: G11 / row for device 526258372598 keeps the old text).
$updatedRows = @(2,3,4,6,7,9,10,12,13)
foreach ($r in $updatedRows) {
    $ws.Cells.Item($r, 7).Value = "Actualizado *"
}

# New observation notes added in column H for rows 10 and 12.
$ws.Cells.Item(10, 8).Value = "Conexión lenta al GPS"
$ws.Cells.Item(12, 8).Value = "Modificar interruptor de encendido"

# Leave the selection on G9, matching the saved view state.
[void]$ws.Range("G9").Select()
